# Add new vocabulary rows (89-152) to Lesson 15 wordlist
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,('to die', '死ぬ|しぬ')
    ,('death', '死|し')
    ,('desperate', '必死|ひっし')
    ,('the dead', '死者|ししゃ')
    ,('meaning', '意味|いみ')
    ,('to watch out', '注意する|ちゅういする')
    ,('opinion', '意見|いけん')
    ,('to prepare', '用意する|よういする')
    ,('hobby', '趣味|しゅみ')
    ,('interest', '興味|きょうみ')
    ,('soybean paste', '味噌|みそ')
    ,('taste', '味|あじ')
    ,('to watch out', '注意する|ちゅういする')
    ,('to order', '注文する|ちゅうもんする')
    ,('to pour', '注ぐ|そそぐ')
    ,('summer', '夏|なつ')
    ,('summer vacation', '夏休み|なつやすみ')
    ,('early summer', '初夏|しょか')
    ,('fish', '魚|さかな')
    ,('fish shop', '魚屋|さかなや')
    ,('gold fish', '金魚|きんぎょ')
    ,('mermaid', '人魚|にんぎょ')
    ,('fish market', '魚市場|うおいちば')
    ,('temple', 'お寺|おてら')
    ,('Toji Temple', '東寺|とうじ')
    ,('sacred building', '寺院|じいん')
    ,('Zen temple', '禅寺|ぜんでら')
    ,('wide; spacious', '広い|ひろい')
    ,('square; open space', '広場|ひろば')
    ,('Hiroshima', '広島|ひろしま')
    ,('advertisement', '広告|こうこく')
    ,('foot; leg', '足|あし')
    ,('to be sufficient', '足りる|たりる')
    ,('one pair of shoes', '一足|いっそく')
    ,('lack of water', '水不足|みずぶそく')
    ,('bicycle', '自転車|じてんしゃ')
    ,('to drive', '運転する|うんてんする')
    ,('conveyor belt sushi', '回転ずし|かいてんずし')
    ,('to tumble; to fall down', '転ぶ|ころぶ')
    ,('to borrow', '借りる|かりる')
    ,('rented land', '借地|しゃくち')
    ,('debt', '借金|しゃっきん')
    ,('rented house', '借家|しゃくや')
    ,('to run', '走る|はしる')
    ,('hasty writing', '走り書き|はしりがき')
    ,('escape from a prison', '脱走|だっそう')
    ,('place', '場所|ばしょ')
    ,('factory', '工場|こうじょう')
    ,('market', '市場|いちば')
    ,('case', '場合|ばあい')
    ,('parking lot', '駐車場|ちゅうしゃじょう')
    ,('building', '建物|たてもの')
    ,('to build', '建てる|たてる')
    ,('to be built', '建つ|たつ')
    ,('founding a nation', '建国|けんこく')
    ,('underground', '地下|ちか')
    ,('subway', '地下鉄|ちかてつ')
    ,('map', '地図|ちず')
    ,('the earth', '地球|ちきゅう')
    ,('earthquake', '地震|じしん')
    ,('to go through; to pass', '通る|とおる')
    ,('to commute', '通う|かよう')
    ,('going to school', '通学|つうがく')
    ,('going to work', '通勤|つうきん')
)

$startRow = 89
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $pair = $data[$i]
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
}

